$d = $word.ActiveDocument

# Update the date header
$found = $d.Content.Find.Execute("2023-05-02 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-05-03 Wednesday", 2)

# Update each multiplication cell in the table, addressed positionally (row, col)
# to avoid any ambiguity from repeated/overlapping text values.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "60×51=3060"
$t.Cell(1, 2).Range.Text = "57×66=3762"
$t.Cell(1, 3).Range.Text = "68×35=2380"
$t.Cell(1, 4).Range.Text = "46×88=4048"
$t.Cell(1, 5).Range.Text = "65×19=1235"

$t.Cell(2, 1).Range.Text = "54×99=5346"
$t.Cell(2, 2).Range.Text = "10×58=580"
$t.Cell(2, 3).Range.Text = "57×78=4446"
$t.Cell(2, 4).Range.Text = "42×81=3402"
$t.Cell(2, 5).Range.Text = "14×90=1260"

$t.Cell(3, 1).Range.Text = "37×97=3589"
$t.Cell(3, 2).Range.Text = "59×46=2714"
$t.Cell(3, 3).Range.Text = "84×99=8316"
$t.Cell(3, 4).Range.Text = "87×79=6873"
$t.Cell(3, 5).Range.Text = "45×41=1845"

$t.Cell(4, 1).Range.Text = "13×79=1027"
$t.Cell(4, 2).Range.Text = "49×72=3528"
$t.Cell(4, 3).Range.Text = "52×39=2028"
$t.Cell(4, 4).Range.Text = "57×63=3591"
$t.Cell(4, 5).Range.Text = "91×59=5369"

$t.Cell(5, 1).Range.Text = "16×54=864"
$t.Cell(5, 2).Range.Text = "93×35=3255"
$t.Cell(5, 3).Range.Text = "50×31=1550"
$t.Cell(5, 4).Range.Text = "37×86=3182"
$t.Cell(5, 5).Range.Text = "97×10=970"

$t.Cell(6, 1).Range.Text = "81×81=6561"
$t.Cell(6, 2).Range.Text = "61×84=5124"
$t.Cell(6, 3).Range.Text = "27×96=2592"
$t.Cell(6, 4).Range.Text = "59×11=649"
$t.Cell(6, 5).Range.Text = "53×93=4929"

$t.Cell(7, 1).Range.Text = "64×94=6016"
$t.Cell(7, 2).Range.Text = "32×76=2432"
$t.Cell(7, 3).Range.Text = "53×84=4452"
$t.Cell(7, 4).Range.Text = "82×94=7708"
$t.Cell(7, 5).Range.Text = "84×28=2352"

$t.Cell(8, 1).Range.Text = "92×51=4692"
$t.Cell(8, 2).Range.Text = "90×53=4770"
$t.Cell(8, 3).Range.Text = "68×94=6392"
$t.Cell(8, 4).Range.Text = "49×76=3724"
$t.Cell(8, 5).Range.Text = "66×17=1122"

$t.Cell(9, 1).Range.Text = "27×21=567"
$t.Cell(9, 2).Range.Text = "91×73=6643"
$t.Cell(9, 3).Range.Text = "25×81=2025"
$t.Cell(9, 4).Range.Text = "40×18=720"
$t.Cell(9, 5).Range.Text = "19×82=1558"

$t.Cell(10, 1).Range.Text = "21×93=1953"
$t.Cell(10, 2).Range.Text = "50×98=4900"
$t.Cell(10, 3).Range.Text = "94×53=4982"
$t.Cell(10, 4).Range.Text = "89×14=1246"
$t.Cell(10, 5).Range.Text = "98×30=2940"

$t.Cell(11, 1).Range.Text = "35×91=3185"
$t.Cell(11, 2).Range.Text = "17×14=238"
$t.Cell(11, 3).Range.Text = "63×56=3528"
$t.Cell(11, 4).Range.Text = "70×41=2870"
$t.Cell(11, 5).Range.Text = "11×79=869"

$t.Cell(12, 1).Range.Text = "57×71=4047"
$t.Cell(12, 2).Range.Text = "25×33=825"
$t.Cell(12, 3).Range.Text = "30×13=390"
$t.Cell(12, 4).Range.Text = "55×43=2365"
$t.Cell(12, 5).Range.Text = "38×43=1634"

$t.Cell(13, 1).Range.Text = "45×23=1035"
$t.Cell(13, 2).Range.Text = "92×50=4600"
$t.Cell(13, 3).Range.Text = "55×55=3025"
$t.Cell(13, 4).Range.Text = "61×55=3355"
$t.Cell(13, 5).Range.Text = "99×68=6732"

$t.Cell(14, 1).Range.Text = "46×70=3220"
$t.Cell(14, 2).Range.Text = "57×27=1539"
$t.Cell(14, 3).Range.Text = "84×63=5292"
$t.Cell(14, 4).Range.Text = "100×44=4400"
$t.Cell(14, 5).Range.Text = "25×88=2200"

$t.Cell(15, 1).Range.Text = "22×15=330"
$t.Cell(15, 2).Range.Text = "41×47=1927"
$t.Cell(15, 3).Range.Text = "11×75=825"
$t.Cell(15, 4).Range.Text = "61×95=5795"
$t.Cell(15, 5).Range.Text = "42×90=3780"

$t.Cell(16, 1).Range.Text = "66×99=6534"
$t.Cell(16, 2).Range.Text = "13×59=767"
$t.Cell(16, 3).Range.Text = "74×39=2886"
$t.Cell(16, 4).Range.Text = "93×25=2325"
$t.Cell(16, 5).Range.Text = "58×46=2668"

$t.Cell(17, 1).Range.Text = "78×72=5616"
$t.Cell(17, 2).Range.Text = "66×98=6468"
$t.Cell(17, 3).Range.Text = "80×44=3520"
$t.Cell(17, 4).Range.Text = "13×69=897"
$t.Cell(17, 5).Range.Text = "92×42=3864"

$t.Cell(18, 1).Range.Text = "77×54=4158"
$t.Cell(18, 2).Range.Text = "48×63=3024"
$t.Cell(18, 3).Range.Text = "71×95=6745"
$t.Cell(18, 4).Range.Text = "78×86=6708"
$t.Cell(18, 5).Range.Text = "58×24=1392"

$t.Cell(19, 1).Range.Text = "47×48=2256"
$t.Cell(19, 2).Range.Text = "26×53=1378"
$t.Cell(19, 3).Range.Text = "32×66=2112"
$t.Cell(19, 4).Range.Text = "95×32=3040"
$t.Cell(19, 5).Range.Text = "68×91=6188"

$t.Cell(20, 1).Range.Text = "76×95=7220"
$t.Cell(20, 2).Range.Text = "25×85=2125"
$t.Cell(20, 3).Range.Text = "13×73=949"
$t.Cell(20, 4).Range.Text = "24×62=1488"
$t.Cell(20, 5).Range.Text = "35×38=1330"
